$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width changes
# NOTE: the runtime stores column width internally quantized to 1/6-character
# increments, so we pick the ColumnWidth input that makes the stored OOXML
# <col width="..."> attribute land as close as possible to the target
# (A -> 15.7109375, B -> 16.42578125).
$ws.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666

# Cell value changes
$ws.Range("A1").Value = -0.11767664784914444
$ws.Range("B1").Value = 0.11725705189250135
$ws.Range("A2").Value = -0.078550708444662831
$ws.Range("B2").Value = 0.077160555320958224
$ws.Range("A3").Value = -0.027446623346666144
$ws.Range("B3").Value = 0.027049036983433083
$ws.Range("A4").Value = -0.019049037106011468
$ws.Range("B4").Value = 0.018691865942484398
$ws.Range("A5").Value = -0.015691865995913545
$ws.Range("B5").Value = 0.014478376903874768
$ws.Range("A6").Value = 0.000526950290701933
$ws.Range("B6").Value = -0.00072905156864955245
$ws.Range("A7").Value = 0.010729051420158786
$ws.Range("B7").Value = -0.010769406998908604
$ws.Range("A8").Value = 0.020769406852954031
$ws.Range("B8").Value = -0.020834907970484817
$ws.Range("A9").Value = 0.022834907943074523
$ws.Range("B9").Value = -0.022888706809638393
$ws.Range("A10").Value = 0.024888706789756299
$ws.Range("B10").Value = -0.024888822699255186
$ws.Range("A11").Value = 0.027888822666127133
$ws.Range("B11").Value = -0.02789494151234706
$ws.Range("A12").Value = -0.004677902604028894
$ws.Range("B12").Value = 0.004676070010806832
$ws.Range("A13").Value = -0.001176070048936495
$ws.Range("B13").Value = 0.0011748179321440233
$ws.Range("A14").Value = 0.0068251819676712344
$ws.Range("B14").Value = -0.0068367804164015666
$ws.Range("A15").Value = 0.0078367804150039078
$ws.Range("B15").Value = -0.0078587967137506709
$ws.Range("A16").Value = -0.0040458542640600115
$ws.Range("B16").Value = 0.0036645348787889276
$ws.Range("A17").Value = -0.0016645348932300408
$ws.Range("B17").Value = 0.0015112604061275903
$ws.Range("A18").Value = -0.032706343533046578
$ws.Range("B18").Value = 0.032658008500366265
$ws.Range("A19").Value = -0.028658008558448245
$ws.Range("B19").Value = 0.02829780326414788
$ws.Range("A20").Value = -0.0080157681977031103
$ws.Range("B20").Value = 0.0080056753737487441
$ws.Range("A21").Value = -0.0040056754433557273
$ws.Range("B21").Value = 0.003999999929812148
$ws.Range("A22").Value = -0.045713932039880945
$ws.Range("B22").Value = 0.045500770971964144
$ws.Range("A23").Value = -0.040500771056323437
$ws.Range("B23").Value = 0.040099116392262601
$ws.Range("A24").Value = -0.020099116693923946
$ws.Range("B24").Value = 0.019999999694007009
$ws.Range("A25").Value = -0.097244879941269957
$ws.Range("B25").Value = 0.097123205860164319
$ws.Range("A26").Value = -0.094623205909698527
$ws.Range("B26").Value = 0.094464521519245892
$ws.Range("A27").Value = -0.091964521575181646
$ws.Range("B27").Value = 0.091011683846474511
$ws.Range("A28").Value = -0.08901168392335812
$ws.Range("B28").Value = 0.088359714789656429
$ws.Range("A29").Value = -0.081359714955101303
$ws.Range("B29").Value = 0.081171558057771875
$ws.Range("A30").Value = -0.021171558971172288
$ws.Range("B30").Value = 0.021021417429853528
$ws.Range("A31").Value = -0.014021417609539455
$ws.Range("B31").Value = 0.014000623002409895
$ws.Range("A32").Value = -0.004000623224566624
$ws.Range("B32").Value = 0.0039999998614757004
